$d = $word.ActiveDocument
$wNS = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Set-ParaXml($index, $inner) {
    $p = $d.Paragraphs.Item($index)
    $xml = "<w:p $wNS>" + $inner + "</w:p>"
    $p.Range.InsertXML($xml)
}

# --- 1) Split author byline into two runs with a proofErr spellStart/spellEnd around "Rahm" ---
Set-ParaXml 1 "<w:r><w:t xml:space='preserve'>Nicolas C. Broeking &amp; Josh </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>Rahm</w:t></w:r><w:proofErr w:type='spellEnd'/>"

# --- Requirement list items 1-12: split out "Must" with gramStart/gramEnd proofErr wraps ---

# Item 1
Set-ParaXml 8 "<w:r><w:t xml:space='preserve'>1.) </w:t></w:r><w:proofErr w:type='gramStart'/><w:r><w:t>Must</w:t></w:r><w:proofErr w:type='gramEnd'/><w:r><w:t xml:space='preserve'> be supported on at least Mac and Linux.</w:t></w:r>"

# Item 2 (also splits "." and "(" into separate runs wrapped with gramStart/gramEnd)
Set-ParaXml 9 "<w:r><w:t>2.) Must allow for up to 2 consecutive players at a time</w:t></w:r><w:proofErr w:type='gramStart'/><w:r><w:t>.</w:t></w:r><w:r><w:t>(</w:t></w:r><w:proofErr w:type='gramEnd'/><w:r><w:t xml:space='preserve'> More players is a stretch  goal)</w:t></w:r>"

# Item 3 (keep leading "3" and ".) " runs as-is, split "Must allow users to start a game.")
Set-ParaXml 10 "<w:r><w:t>3</w:t></w:r><w:r><w:t xml:space='preserve'>.) </w:t></w:r><w:proofErr w:type='gramStart'/><w:r><w:t>Must</w:t></w:r><w:proofErr w:type='gramEnd'/><w:r><w:t xml:space='preserve'> allow users to start a game.</w:t></w:r>"

# Item 4
Set-ParaXml 11 "<w:r><w:t xml:space='preserve'>4.) </w:t></w:r><w:proofErr w:type='gramStart'/><w:r><w:t>Must</w:t></w:r><w:proofErr w:type='gramEnd'/><w:r><w:t xml:space='preserve'> allow users to join a game.</w:t></w:r>"

# Item 5 (keep "game state" and "." runs)
Set-ParaXml 12 "<w:r><w:t xml:space='preserve'>5.) </w:t></w:r><w:proofErr w:type='gramStart'/><w:r><w:t>Must</w:t></w:r><w:proofErr w:type='gramEnd'/><w:r><w:t xml:space='preserve'> allow users to interact with </w:t></w:r><w:r><w:t>game state</w:t></w:r><w:r><w:t>.</w:t></w:r>"

# Item 6
Set-ParaXml 13 "<w:r><w:t xml:space='preserve'>6.) </w:t></w:r><w:proofErr w:type='gramStart'/><w:r><w:t>Must</w:t></w:r><w:proofErr w:type='gramEnd'/><w:r><w:t xml:space='preserve'> allow users to place troops.</w:t></w:r>"

# Item 7 (keep "countries" and "." runs)
Set-ParaXml 14 "<w:r><w:t xml:space='preserve'>7.) </w:t></w:r><w:proofErr w:type='gramStart'/><w:r><w:t>Must</w:t></w:r><w:proofErr w:type='gramEnd'/><w:r><w:t xml:space='preserve'> allow users to attack </w:t></w:r><w:r><w:t>countries</w:t></w:r><w:r><w:t>.</w:t></w:r>"

# Item 8
Set-ParaXml 15 "<w:r><w:t xml:space='preserve'>8.) </w:t></w:r><w:proofErr w:type='gramStart'/><w:r><w:t>Must</w:t></w:r><w:proofErr w:type='gramEnd'/><w:r><w:t xml:space='preserve'> allow users to end attack phase.</w:t></w:r>"

# Item 9
Set-ParaXml 16 "<w:r><w:t xml:space='preserve'>9.) </w:t></w:r><w:proofErr w:type='gramStart'/><w:r><w:t>Must</w:t></w:r><w:proofErr w:type='gramEnd'/><w:r><w:t xml:space='preserve'> allow users to move troops.</w:t></w:r>"

# Item 10
Set-ParaXml 17 "<w:r><w:t xml:space='preserve'>10.) </w:t></w:r><w:proofErr w:type='gramStart'/><w:r><w:t>Must</w:t></w:r><w:proofErr w:type='gramEnd'/><w:r><w:t xml:space='preserve'> allow users to end move phase.</w:t></w:r>"

# Item 11 (keep "game state" and "." runs)
Set-ParaXml 18 "<w:r><w:t xml:space='preserve'>11.) </w:t></w:r><w:proofErr w:type='gramStart'/><w:r><w:t>Must</w:t></w:r><w:proofErr w:type='gramEnd'/><w:r><w:t xml:space='preserve'> have the system manage the </w:t></w:r><w:r><w:t>game state</w:t></w:r><w:r><w:t>.</w:t></w:r>"

# Item 12
Set-ParaXml 19 "<w:r><w:t xml:space='preserve'>12.) </w:t></w:r><w:proofErr w:type='gramStart'/><w:r><w:t>Must</w:t></w:r><w:proofErr w:type='gramEnd'/><w:r><w:t xml:space='preserve'> have the system organize turns. </w:t></w:r>"

# --- Architecture paragraph: split out "openGL" with spellStart/spellEnd ---
Set-ParaXml 53 "<w:r><w:t xml:space='preserve'>We will have two kinds of systems, Server systems and client systems. The servers and clients will communicate by sending XML through sockets to each other. The client side will use SDL and </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>openGL</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t xml:space='preserve'> to display the graphics and run the event loop.</w:t></w:r>"

# --- Final section: extend "VI." paragraph and append the new Data Storage /
#     UI Mockups / UI Interactions sections, re-appending the _GoBack bookmark
#     in its own trailing paragraph at the very end. ---
$viIndex = $d.Paragraphs.Count
$viXml = @"
<w:p $wNS>
  <w:r><w:t>VI.</w:t></w:r>
  <w:r><w:t xml:space='preserve'> Data Storage</w:t></w:r>
</w:p>
<w:p $wNS/>
<w:p $wNS>
  <w:r>
    <w:tab/>
    <w:t>We do not need to save data in any way but we do need to send information between the server and the client. We are going to do this using XML.</w:t>
  </w:r>
</w:p>
<w:p $wNS/>
<w:p $wNS>
  <w:r><w:t xml:space='preserve'>VII. </w:t></w:r>
  <w:r><w:t>UI Mockups</w:t></w:r>
</w:p>
<w:p $wNS/>
<w:p $wNS>
  <w:r>
    <w:tab/>
    <w:t xml:space='preserve'>There are two kinds of views for risk. You can either play risk using a </w:t>
  </w:r>
  <w:proofErr w:type='spellStart'/>
  <w:proofErr w:type='gramStart'/>
  <w:r><w:t>gui</w:t></w:r>
  <w:proofErr w:type='spellEnd'/>
  <w:proofErr w:type='gramEnd'/>
  <w:r><w:t xml:space='preserve'> that opens up to the game with a </w:t></w:r>
  <w:proofErr w:type='spellStart'/>
  <w:r><w:t>hud</w:t></w:r>
  <w:proofErr w:type='spellEnd'/>
  <w:r><w:t xml:space='preserve'> located at the bottom or you can play using the command line. </w:t></w:r>
</w:p>
<w:p $wNS/>
<w:p $wNS>
  <w:r><w:t>The UI Mockup is in the UI Mockups Folder.</w:t></w:r>
</w:p>
<w:p $wNS/>
<w:p $wNS>
  <w:r><w:t>VIII.</w:t></w:r>
  <w:r><w:t xml:space='preserve'> UI Interactions</w:t></w:r>
</w:p>
<w:p $wNS/>
<w:p $wNS>
  <w:r>
    <w:tab/>
    <w:t xml:space='preserve'>The user interacts with the system in one way. They must click on a region to attack. When clicked the event loop sends a signal to the mouse handler object. This is shown in the User Interaction Sequence Diagram found in the Diagrams folder. Then the system must figure out what the user is trying to do and update the game state. This is shown in the update server sequence diagram. Then the server must validate the move and then update the game state. This is shown in the server validation sequence diagram. </w:t>
  </w:r>
</w:p>
<w:p $wNS/>
<w:p $wNS>
  <w:bookmarkStart w:id='0' w:name='_GoBack'/>
  <w:bookmarkEnd w:id='0'/>
</w:p>
"@
$d.Paragraphs.Item($viIndex).Range.InsertXML($viXml)

Write-Host "done steps 1-15"
